$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last refreshed" timestamp shown in row 1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 11:22"

# Refresh case counts and re-sort-driven country reassignments for the rows
# whose rank (by column B, "Casos totales" desc) changed between the 10:52 and
# 11:22 snapshots. Each row below keeps its position but gets the country name
# and B:H stats that now belong at that rank.
$ws.Cells.Item(13, 1).Value = "Belgica"
$ws.Cells.Item(13, 2).Value = 30589
$ws.Cells.Item(13, 3).Value = 942
$ws.Cells.Item(13, 4).Value = 6707
$ws.Cells.Item(13, 5).Value = 19979
$ws.Cells.Item(13, 6).Value = 1234
$ws.Cells.Item(13, 7).Value = 303
$ws.Cells.Item(13, 8).Value = 3903

$ws.Cells.Item(15, 1).Value = "Suiza"
$ws.Cells.Item(15, 2).Value = 25449
$ws.Cells.Item(15, 3).Value = 34
$ws.Cells.Item(15, 4).Value = 12700
$ws.Cells.Item(15, 5).Value = 11634
$ws.Cells.Item(15, 6).Value = 386
$ws.Cells.Item(15, 7).Value = 9
$ws.Cells.Item(15, 8).Value = 1115

$ws.Cells.Item(20, 1).Value = "Austria"
$ws.Cells.Item(20, 2).Value = 13972
$ws.Cells.Item(20, 3).Value = 27
$ws.Cells.Item(20, 4).Value = 7343
$ws.Cells.Item(20, 5).Value = 6261
$ws.Cells.Item(20, 6).Value = 239
$ws.Cells.Item(20, 7).Value = 18
$ws.Cells.Item(20, 8).Value = 368

$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 4557
$ws.Cells.Item(40, 3).Value = 316
$ws.Cells.Item(40, 4).Value = 380
$ws.Cells.Item(40, 5).Value = 3778
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 26
$ws.Cells.Item(40, 8).Value = 399

$ws.Cells.Item(41, 1).Value = "Arabia Saudita"
$ws.Cells.Item(41, 2).Value = 4462
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 761
$ws.Cells.Item(41, 5).Value = 3642
$ws.Cells.Item(41, 6).Value = 67
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 59

$ws.Cells.Item(47, 1).Value = "Finlandia"
$ws.Cells.Item(47, 2).Value = 3064
$ws.Cells.Item(47, 3).Value = 90
$ws.Cells.Item(47, 4).Value = 300
$ws.Cells.Item(47, 5).Value = 2708
$ws.Cells.Item(47, 6).Value = 80
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 56

$ws.Cells.Item(48, 1).Value = "Catar"
$ws.Cells.Item(48, 2).Value = 2979
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 275
$ws.Cells.Item(48, 5).Value = 2697
$ws.Cells.Item(48, 6).Value = 37
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 7

$ws.Cells.Item(59, 1).Value = "Marruecos"
$ws.Cells.Item(59, 2).Value = 1746
$ws.Cells.Item(59, 3).Value = 85
$ws.Cells.Item(59, 4).Value = 196
$ws.Cells.Item(59, 5).Value = 1430
$ws.Cells.Item(59, 6).Value = 1
$ws.Cells.Item(59, 7).Value = 2
$ws.Cells.Item(59, 8).Value = 120

$ws.Cells.Item(60, 1).Value = "Islandia"
$ws.Cells.Item(60, 2).Value = 1701
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 889
$ws.Cells.Item(60, 5).Value = 804
$ws.Cells.Item(60, 6).Value = 10
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 8

$ws.Cells.Item(61, 1).Value = "Moldavia"
$ws.Cells.Item(61, 2).Value = 1662
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 107
$ws.Cells.Item(61, 5).Value = 1522
$ws.Cells.Item(61, 6).Value = 80
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 33

$ws.Cells.Item(114, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(114, 2).Value = 235
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 17
$ws.Cells.Item(114, 5).Value = 198
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 20

$ws.Cells.Item(166, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(166, 2).Value = 23
$ws.Cells.Item(166, 3).Value = 2
$ws.Cells.Item(166, 4).Value = 0
$ws.Cells.Item(166, 5).Value = 21
$ws.Cells.Item(166, 6).Value = 1
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 2

$ws.Cells.Item(167, 1).Value = "Mozambique"
$ws.Cells.Item(167, 2).Value = 21
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = 19
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(193, 2).Value = 11
$ws.Cells.Item(193, 3).Value = 3
$ws.Cells.Item(193, 4).Value = 3
$ws.Cells.Item(193, 5).Value = 8
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "Groenlandia"
$ws.Cells.Item(194, 2).Value = 11
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 11
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(195, 1).Value = "Sierra Leona"
$ws.Cells.Item(195, 2).Value = 10
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

$ws.Cells.Item(196, 1).Value = "Surinam"
$ws.Cells.Item(196, 2).Value = 10
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 4
$ws.Cells.Item(196, 5).Value = 5
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1

$ws.Cells.Item(197, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(197, 2).Value = 9
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 1

$ws.Cells.Item(198, 1).Value = "Gambia"
$ws.Cells.Item(198, 2).Value = 9
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 2
$ws.Cells.Item(198, 5).Value = 6
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(199, 1).Value = "Nicaragua"
$ws.Cells.Item(199, 2).Value = 9
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 4
$ws.Cells.Item(199, 5).Value = 4
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 1

$ws.Cells.Item(210, 1).Value = "Timor Oriental"
$ws.Cells.Item(210, 2).Value = 4
$ws.Cells.Item(210, 3).Value = 2
$ws.Cells.Item(210, 4).Value = 1
$ws.Cells.Item(210, 5).Value = 3
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(211, 2).Value = 3
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 0
$ws.Cells.Item(211, 5).Value = 3
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

$ws.Cells.Item(212, 1).Value = "Anguila"
$ws.Cells.Item(212, 2).Value = 3
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 1
$ws.Cells.Item(212, 5).Value = 2
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 0

$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 2).Value = 3
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 2
$ws.Cells.Item(213, 5).Value = 1
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 2).Value = 2
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 0
$ws.Cells.Item(214, 5).Value = 2
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

